$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values in B2:E2
$ws.Range("B2").Value = 10.787023792862547
$ws.Range("C2").Value = 15.393399526523226
$ws.Range("D2").Value = 8.7728464407593894
$ws.Range("E2").Value = 14.269115938797395

# Update row 3 values in B3:E3
$ws.Range("B3").Value = 13.299741894026944
$ws.Range("C3").Value = 15.880770946837549
$ws.Range("D3").Value = 15.748702774842679
$ws.Range("E3").Value = 15.504003619947063

# Update the selection to reflect the new active range B1:E3
$ws.Range("B1:E3").Select()
